$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New historical fy_year rows pushing dates back to 1997-98.
# Rows 17-18 carry the same lito figures as the other "0.04 taper" years,
# rows 19-21 only have the fy_year filled in (values unknown/blank).

$ws.Range("A17").Value = "2001-02"
$ws.Range("B17").Value = 150
$ws.Range("C17").Value = 0.04
$ws.Range("D17").Value = 20700

$ws.Range("A18").Value = "2000-01"
$ws.Range("B18").Value = 150
$ws.Range("C18").Value = 0.04
$ws.Range("D18").Value = 20700

$ws.Range("A19").Value = "1999-00"
$ws.Range("A20").Value = "1998-99"
$ws.Range("A21").Value = "1997-98"

$ws.Range("A22").Select()
